# auto: formatting changes and removing some reports
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 47 (begin group "call"): appearance changes to "field-list summary"
$ws.Range("E47").Value = "field-list summary"

# --- Row 48 (note "first_name"): shorten label, drop appearance value
$ws.Range("C48").Value = ' <i style=”background-color: yellow;”>**${fast_name_ctx}**</i>'
$ws.Range("E48").Value = ""

# --- Insert a brand new row 53 (a highlighted "special instruction" title/separator row)
$ws.Rows.Item(53).Insert()

$ws.Range("A53").Value = "note"
$ws.Range("B53").Value = "n_special_instruction_title"
$ws.Range("C53").Value = "NO_LABEL"
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = "h1 yellow"

# Formatting for the new row: A53 (muted small gray label), B53/D53/E53 (small white-filled
# cells), C53 (plain), and a wide band of blank, orange-filled cells out to AA53.
$rngA = $ws.Range("A53")
$rngA.Font.Name = "Arial"
$rngA.Font.Size = 8
$rngA.Font.Color = 13421772
$rngA.Interior.Color = 13493756
$rngA.WrapText = $true

$rngB = $ws.Range("B53")
$rngB.Font.Name = "Arial"
$rngB.Font.Size = 8
$rngB.Font.ThemeColor = 1
$rngB.Interior.Color = 16777215
$rngB.WrapText = $true

$rngC = $ws.Range("C53")
$rngC.Font.Name = "Arial"
$rngC.Font.ThemeColor = 1

$rngDtoAA = $ws.Range("D53:AA53")
$rngDtoAA.Font.Name = "Arial"
$rngDtoAA.Font.Size = 8
$rngDtoAA.Font.ThemeColor = 1
$rngDtoAA.Interior.Color = 16777215
$rngDtoAA.WrapText = $true

# Data validation dropdown (yes/no) on the relevant cell of the new row
$rngD = $ws.Range("D53")
$rngD.Validation.Add(3, 1, 1, '"yes,no"')
$rngD.Validation.ShowInput = $false
$rngD.Validation.ShowError = $false

Write-Output "edit applied"
